$d = $word.ActiveDocument

$replacements = @(
    @{old="924÷2="; new="825÷4="},
    @{old="889÷2="; new="469÷7="},
    @{old="370÷2="; new="291÷5="},
    @{old="977÷7="; new="656÷8="},
    @{old="470÷6="; new="776÷7="},
    @{old="460÷3="; new="252÷9="},
    @{old="111÷8="; new="730÷8="},
    @{old="836÷2="; new="824÷4="},
    @{old="409÷8="; new="185÷3="},
    @{old="452÷8="; new="630÷7="},
    @{old="351÷4="; new="960÷4="},
    @{old="741÷6="; new="696÷6="},
    @{old="381÷9="; new="374÷7="},
    @{old="296÷9="; new="357÷6="},
    @{old="669÷2="; new="878÷6="},
    @{old="531÷2="; new="551÷6="},
    @{old="335÷2="; new="883÷7="},
    @{old="815÷8="; new="573÷2="},
    @{old="486÷2="; new="412÷3="},
    @{old="865÷8="; new="819÷3="},
    @{old="381÷6="; new="354÷2="},
    @{old="507÷2="; new="442÷2="},
    @{old="679÷8="; new="335÷8="},
    @{old="472÷5="; new="524÷3="},
    @{old="676÷3="; new="164÷6="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
